$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 32.935331
$ws.Range("H2").Value = 98.805993
$ws.Range("I2").Value = 0.1836164637112342
$ws.Range("J2").Value = 0.1836164637112342
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 27.81717833333333
$ws.Range("N2").Value = 83.45153500000001
$ws.Range("O2").Value = 0.4044740580248731
$ws.Range("P2").Value = 0.4044740580248732
$ws.Range("Q2").Value = 916.1679758943617
$ws.Range("R2").Value = 8245.511783049256
$ws.Range("S2").Value = 0.07426809619745973
$ws.Range("T2").Value = 0.07426809619745975

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 32.935331
$ws.Range("H3").Value = 98.805993
$ws.Range("I3").Value = 0.1836164637112342
$ws.Range("J3").Value = 0.1836164637112342
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 4.827410666666666
$ws.Range("N3").Value = 14.482232
$ws.Range("O3").Value = 0.0701926830500802
$ws.Range("P3").Value = 0.0701926830500802
$ws.Range("Q3").Value = 158.9923681795973
$ws.Range("R3").Value = 1430.931313616376
$ws.Range("S3").Value = 0.01288853224005921
$ws.Range("T3").Value = 0.01288853224005921

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 32.935331
$ws.Range("H4").Value = 98.805993
$ws.Range("I4").Value = 0.1836164637112342
$ws.Range("J4").Value = 0.1836164637112342
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 32.06242466666667
$ws.Range("N4").Value = 96.187274
$ws.Range("O4").Value = 0.4662018145637509
$ws.Range("P4").Value = 0.466201814563751
$ws.Range("Q4").Value = 1055.986569059231
$ws.Range("R4").Value = 9503.879121533082
$ws.Range("S4").Value = 0.0856023285659565
$ws.Range("T4").Value = 0.0856023285659565

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 32.935331
$ws.Range("H5").Value = 98.805993
$ws.Range("I5").Value = 0.1836164637112342
$ws.Range("J5").Value = 0.1836164637112342
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.066688333333334
$ws.Range("N5").Value = 12.200065
$ws.Range("O5").Value = 0.05913144436129575
$ws.Range("P5").Value = 0.05913144436129575
$ws.Range("Q5").Value = 133.9377263321717
$ws.Range("R5").Value = 1205.439536989545
$ws.Range("S5").Value = 0.01085750670775872
$ws.Range("T5").Value = 0.01085750670775872

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 95.562134
$ws.Range("H6").Value = 286.686402
$ws.Range("I6").Value = 0.5327646808765668
$ws.Range("J6").Value = 0.5327646808765667
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 27.81717833333333
$ws.Range("N6").Value = 83.45153500000001
$ws.Range("O6").Value = 0.4044740580248731
$ws.Range("P6").Value = 0.4044740580248732
$ws.Range("Q6").Value = 2658.268923391897
$ws.Range("R6").Value = 23924.42031052707
$ws.Range("S6").Value = 0.2154894924464715
$ws.Range("T6").Value = 0.2154894924464715

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 95.562134
$ws.Range("H7").Value = 286.686402
$ws.Range("I7").Value = 0.5327646808765668
$ws.Range("J7").Value = 0.5327646808765667
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.827410666666666
$ws.Range("N7").Value = 14.482232
$ws.Range("O7").Value = 0.0701926830500802
$ws.Range("P7").Value = 0.0701926830500802
$ws.Range("Q7").Value = 461.3176650010293
$ws.Range("R7").Value = 4151.858985009264
$ws.Range("S7").Value = 0.03739618238504597
$ws.Range("T7").Value = 0.03739618238504597

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 95.562134
$ws.Range("H8").Value = 286.686402
$ws.Range("I8").Value = 0.5327646808765668
$ws.Range("J8").Value = 0.5327646808765667
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 32.06242466666667
$ws.Range("N8").Value = 96.187274
$ws.Range("O8").Value = 0.4662018145637509
$ws.Range("P8").Value = 0.466201814563751
$ws.Range("Q8").Value = 3063.953722360905
$ws.Range("R8").Value = 27575.58350124815
$ws.Range("S8").Value = 0.2483758609601331
$ws.Range("T8").Value = 0.2483758609601331

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 95.562134
$ws.Range("H9").Value = 286.686402
$ws.Range("I9").Value = 0.5327646808765668
$ws.Range("J9").Value = 0.5327646808765667
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.066688333333334
$ws.Range("N9").Value = 12.200065
$ws.Range("O9").Value = 0.05913144436129575
$ws.Range("P9").Value = 0.05913144436129575
$ws.Range("Q9").Value = 388.6214154462368
$ws.Range("R9").Value = 3497.592739016131
$ws.Range("S9").Value = 0.0315031450849162
$ws.Range("T9").Value = 0.03150314508491619

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 37.06916133333333
$ws.Range("H10").Value = 111.207484
$ws.Range("I10").Value = 0.2066628180165514
$ws.Range("J10").Value = 0.2066628180165514
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 27.81717833333333
$ws.Range("N10").Value = 83.45153500000001
$ws.Range("O10").Value = 0.4044740580248731
$ws.Range("P10").Value = 0.4044740580248732
$ws.Range("Q10").Value = 1031.159471476438
$ws.Range("R10").Value = 9280.43524328794
$ws.Range("S10").Value = 0.08358974864601042
$ws.Range("T10").Value = 0.08358974864601042

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 37.06916133333333
$ws.Range("H11").Value = 111.207484
$ws.Range("I11").Value = 0.2066628180165514
$ws.Range("J11").Value = 0.2066628180165514
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 4.827410666666666
$ws.Range("N11").Value = 14.482232
$ws.Range("O11").Value = 0.0701926830500802
$ws.Range("P11").Value = 0.0701926830500802
$ws.Range("Q11").Value = 178.9480648249209
$ws.Range("R11").Value = 1610.532583424288
$ws.Range("S11").Value = 0.0145062176832722
$ws.Range("T11").Value = 0.0145062176832722

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 37.06916133333333
$ws.Range("H12").Value = 111.207484
$ws.Range("I12").Value = 0.2066628180165514
$ws.Range("J12").Value = 0.2066628180165514
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 32.06242466666667
$ws.Range("N12").Value = 96.187274
$ws.Range("O12").Value = 0.4662018145637509
$ws.Range("P12").Value = 0.466201814563751
$ws.Range("Q12").Value = 1188.527192706513
$ws.Range("R12").Value = 10696.74473435862
$ws.Range("S12").Value = 0.09634658076217451
$ws.Range("T12").Value = 0.09634658076217451

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 37.06916133333333
$ws.Range("H13").Value = 111.207484
$ws.Range("I13").Value = 0.2066628180165514
$ws.Range("J13").Value = 0.2066628180165514
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 4.066688333333334
$ws.Range("N13").Value = 12.200065
$ws.Range("O13").Value = 0.05913144436129575
$ws.Range("P13").Value = 0.05913144436129575
$ws.Range("Q13").Value = 150.7487259207178
$ws.Range("R13").Value = 1356.73853328646
$ws.Range("S13").Value = 0.0122202709250943
$ws.Range("T13").Value = 0.0122202709250943

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 13.80362366666667
$ws.Range("H14").Value = 41.410871
$ws.Range("I14").Value = 0.07695603739564764
$ws.Range("J14").Value = 0.07695603739564763
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 27.81717833333333
$ws.Range("N14").Value = 83.45153500000001
$ws.Range("O14").Value = 0.4044740580248731
$ws.Range("P14").Value = 0.4044740580248732
$ws.Range("Q14").Value = 383.9778611818872
$ws.Range("R14").Value = 3455.800750636985
$ws.Range("S14").Value = 0.03112672073493149
$ws.Range("T14").Value = 0.03112672073493149

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 13.80362366666667
$ws.Range("H15").Value = 41.410871
$ws.Range("I15").Value = 0.07695603739564764
$ws.Range("J15").Value = 0.07695603739564763
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 4.827410666666666
$ws.Range("N15").Value = 14.482232
$ws.Range("O15").Value = 0.0701926830500802
$ws.Range("P15").Value = 0.0701926830500802
$ws.Range("Q15").Value = 66.63576012711911
$ws.Range("R15").Value = 599.721841144072
$ws.Range("S15").Value = 0.005401750741702814
$ws.Range("T15").Value = 0.005401750741702813

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 13.80362366666667
$ws.Range("H16").Value = 41.410871
$ws.Range("I16").Value = 0.07695603739564764
$ws.Range("J16").Value = 0.07695603739564763
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 32.06242466666667
$ws.Range("N16").Value = 96.187274
$ws.Range("O16").Value = 0.4662018145637509
$ws.Range("P16").Value = 0.466201814563751
$ws.Range("Q16").Value = 442.5776439395171
$ws.Range("R16").Value = 3983.198795455654
$ws.Range("S16").Value = 0.0358770442754868
$ws.Range("T16").Value = 0.0358770442754868

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 13.80362366666667
$ws.Range("H17").Value = 41.410871
$ws.Range("I17").Value = 0.07695603739564764
$ws.Range("J17").Value = 0.07695603739564763
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 4.066688333333334
$ws.Range("N17").Value = 12.200065
$ws.Range("O17").Value = 0.05913144436129575
$ws.Range("P17").Value = 0.05913144436129575
$ws.Range("Q17").Value = 56.13503532295724
$ws.Range("R17").Value = 505.2153179066151
$ws.Range("S17").Value = 0.004550521643526533
$ws.Range("T17").Value = 0.004550521643526532

